$wb = $excel.ActiveWorkbook
$storeWs = $wb.Worksheets.Item("store")
$menuWs = $wb.Worksheets.Item("menu")

# Replace the bare "xxx.jpg" filenames in the store sheet's image column (C)
# with real hyperlinked image URLs, matching what the "menu" sheet already does
# for its own image column.
$imageUrls = @{
    2 = "https://sebenarnya.my/wp-content/uploads/2017/02/old-town-coffee-1.jpg"
    3 = "https://media-cdn.tripadvisor.com/media/photo-s/06/a4/b3/f2/kfc-kentucky-fried-chicken.jpg"
    4 = "https://media-cdn.tripadvisor.com/media/photo-s/1a/9f/86/0f/caption.jpg"
    5 = "https://www.thedailymeal.com/img/gallery/avoid-these-14-mistakes-when-ordering-at-mcdonalds/intro-1680098756.jpg"
    6 = "https://proriat-franchise.com/wp-content/uploads/2021/08/161632.jpg"
    7 = "https://media-cdn.tripadvisor.com/media/photo-s/0e/47/81/3e/outside.jpg"
}

foreach ($row in 2..7) {
    $cell = $storeWs.Range("C$row")
    $url = $imageUrls[$row]
    $cell.Value = $url
    $storeWs.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}

# Resize the store sheet's columns to fit the new, much wider URL column.
$storeWs.Columns.Item(2).ColumnWidth = 18.8
$storeWs.Columns.Item(3).ColumnWidth = 97.15
$storeWs.Columns.Item(4).ColumnWidth = 9.3

# The "My Account" page now reads the store sheet instead of the menu sheet,
# so make "store" the active/selected sheet and "menu" inactive.
$storeWs.Activate()
$storeWs.Range("C11").Select()
